$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.144.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.56%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.653.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.44%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'214.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.40%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5108"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.66%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.17%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.26%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -4.26%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.57%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.56%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.648.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.78%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.277"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -5.29%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.882.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.36%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.5505"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -5.68%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅7999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.75%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'63.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -6.29%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'26.176.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.39%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.17%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'208.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -6.13%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.409"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.01%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.68%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.00%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.05%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'143.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.17%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +2.34%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.1166"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.79%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.962"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.93%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.05131"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.47%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -4.22%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.23%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.213"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.46%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.95%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -4.13%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.374"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.87%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.9273"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.69%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5685"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.04%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.157.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +5.85%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01589"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.97%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.05%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.8338"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.27%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.640"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.99%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'100.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.93%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.791.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.42%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.67%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.4548"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.33%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'55.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'1.005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.10%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'7.856"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.86%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.05038"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.77%  "
$ws.Range("E51").Style = "Normal"
